$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 511, shifting existing rows 511:551 down to 512:552.
$ws.Rows("511:511").Insert()

# Populate the newly inserted row 511 with its data.
$ws.Range("A511").Value = 6
$ws.Range("B511").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C511").Value = 'Metropolitana'
$ws.Range("D511").Value = 45013
$ws.Range("E511").Value = 13
$ws.Range("F511").Value = 100112032
$ws.Range("G511").Value = 'Zapallo italiano'
$ws.Range("H511").Value = 'Sin especificar'
$ws.Range("I511").Value = 'Primera'
$ws.Range("J511").Value = 610
$ws.Range("K511").Value = 5000
$ws.Range("L511").Value = 6000
$ws.Range("M511").Value = 5475
$ws.Range("N511").Value = '$/caja 50 unidades'
$ws.Range("O511").Value = 'Región Metropolitana'
$ws.Range("P511").Value = 110
$ws.Range("Q511").Value = 50
$ws.Range("R511").Value = 'Hortaliza'
